$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.180.61"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.577.05"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.40"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.72"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("E8").Value = "  +2.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.77"
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0993"
$ws.Range("E10").Value = "  -4.05%  "
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.331"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.032.81"
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "58.110.09"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.62"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.584.05"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  -2.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.40"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "334.39"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.01"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.15"
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.62"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.417"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.995"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("E26").Value = "  -4.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.04"
$ws.Range("E27").Value = "  -3.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0723"
$ws.Range("E29").Value = "  -4.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.63"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.37"
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.83"
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.84"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.87"
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "36.86"
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("E36").Value = "  -4.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.827"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.816"
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.58"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "283.43"
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.588"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.64"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0945"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0532"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.31"
$ws.Range("E47").Value = "  -4.52%  "
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.903.67"
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.74"
$ws.Range("E50").Value = "  -3.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.36"
$ws.Range("E51").Value = "  -4.40%  "

Write-Host "Applied cryptos update"
